# Flash Header.xlsx - add command table entries (Command / Transmit File BIN /
# Transmit File ASCII / Wipe disk) pulled in from the engine controller code base,
# and switch the active sheet/selection over to the "Serial Command Encoding" tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$wsCmd = $wb.Worksheets.Item("Serial Command Encoding")

# --- Update the "Serial Command Encoding" table ------------------------------
# New strings are introduced in this order so the shared-string table ends up
# ordered the same way as in the authored workbook: Command, Transmit File BIN,
# Transmit File ASCII, Wipe disk, 0xff.

# Header: "Encoding" -> "Command"
$wsCmd.Range("A1").Value = "Command"

# Row 3 used to be the single "Transmit File" command; it becomes
# "Transmit File BIN".
$wsCmd.Range("A3").Value = "Transmit File BIN"

# New row 4: "Transmit File ASCII", byte 2, same File Number arg.
$wsCmd.Range("A4").Value = "Transmit File ASCII"
$wsCmd.Range("B4").Value = 2
$wsCmd.Range("C4").Value = "File Number (uint8)"

# New row 5: "Wipe disk", byte 3.
$wsCmd.Range("A5").Value = "Wipe disk"
$wsCmd.Range("B5").Value = 3

# 0xff terminator byte, now documented on rows 2, 3, 4 and 5.
$wsCmd.Range("C2").Value = "0xff"
$wsCmd.Range("D3").Value = "0xff"
$wsCmd.Range("D4").Value = "0xff"
$wsCmd.Range("C5").Value = "0xff"

# --- Switch the active tab / selection ---------------------------------------
# Previously "Sheet1" was the selected tab with A2:A7 selected; now the
# "Serial Command Encoding" tab is selected (activeTab=1) with C9 selected,
# and Sheet1 keeps a plain B9 selection.
$ws1.Activate()
$ws1.Range("B9").Select() | Out-Null

$wsCmd.Activate()
$wsCmd.Range("C9").Select() | Out-Null
